# Applies the per-cell numeric corrections from the CryCompanywiseStockReport_1
# revision: quantity (F) / amount (G) restatements on individual stock-item rows,
# the resulting Sub Total (B) / Grand Total (B) roll-ups, and a handful of rows
# (e.g. 283/284, 448/449, 456/457, 563/564, 573/574, 721/722) whose two stock-code
# entries were swapped (code/rate/qty/amount interchanged between the pair).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: F6=35, G6=1045.8
$ws.Range("F6").Value = 35
$ws.Range("G6").Value = 1045.8

# Row 10: B10=24989.32
$ws.Range("B10").Value = 24989.32

# Row 89: F89=92, G89=13861.64
$ws.Range("F89").Value = 92
$ws.Range("G89").Value = 13861.64

# Row 90: F90=1, G90=102.46
$ws.Range("F90").Value = 1
$ws.Range("G90").Value = 102.46

# Row 96: B96=154474.95
$ws.Range("B96").Value = 154474.95

# Row 113: F113=6, G113=8463.9
$ws.Range("F113").Value = 6
$ws.Range("G113").Value = 8463.9

# Row 115: B115=77027.14999999999
$ws.Range("B115").Value = 77027.14999999999

# Row 119: F119=224, G119=11661.44
$ws.Range("F119").Value = 224
$ws.Range("G119").Value = 11661.44

# Row 124: F124=53, G124=3480.51
$ws.Range("F124").Value = 53
$ws.Range("G124").Value = 3480.51

# Row 125: B125=55249.02
$ws.Range("B125").Value = 55249.02

# Row 140: F140=156, G140=15102.36
$ws.Range("F140").Value = 156
$ws.Range("G140").Value = 15102.36

# Row 142: B142=8943.42
$ws.Range("B142").Value = 8943.42

# Row 185: F185=330, G185=2636.7
$ws.Range("F185").Value = 330
$ws.Range("G185").Value = 2636.7

# Row 187: B187=10745.29
$ws.Range("B187").Value = 10745.29

# Row 189: F189=187, G189=12117.6
$ws.Range("F189").Value = 187
$ws.Range("G189").Value = 12117.6

# Row 196: B196=25352.14
$ws.Range("B196").Value = 25352.14

# Row 283: B283=55373, E283=163.62, F283=-94, G283=-13562.32
$ws.Range("B283").Value = 55373
$ws.Range("E283").Value = 163.62
$ws.Range("F283").Value = -94
$ws.Range("G283").Value = -13562.32

# Row 284: B284=63520, E284=153.4, F284=63, G284=9089.639999999999
$ws.Range("B284").Value = 63520
$ws.Range("E284").Value = 153.4
$ws.Range("F284").Value = 63
$ws.Range("G284").Value = 9089.639999999999

# Row 352: F352=29, G352=2414.83
$ws.Range("F352").Value = 29
$ws.Range("G352").Value = 2414.83

# Row 354: F354=11, G354=784.96
$ws.Range("F354").Value = 11
$ws.Range("G354").Value = 784.96

# Row 364: B364=149665.39
$ws.Range("B364").Value = 149665.39

# Row 415: F415=16, G415=982.5599999999999
$ws.Range("F415").Value = 16
$ws.Range("G415").Value = 982.5599999999999

# Row 416: B416=21323.81
$ws.Range("B416").Value = 21323.81

# Row 433: F433=1, G433=2549.35
$ws.Range("F433").Value = 1
$ws.Range("G433").Value = 2549.35

# Row 444: F444=1, G444=3484.82
$ws.Range("F444").Value = 1
$ws.Range("G444").Value = 3484.82

# Row 445: B445=241720.77
$ws.Range("B445").Value = 241720.77

# Row 448: B448=53602, E448=15.69, F448=-231, G448=-3037.65
$ws.Range("B448").Value = 53602
$ws.Range("E448").Value = 15.69
$ws.Range("F448").Value = -231
$ws.Range("G448").Value = -3037.65

# Row 449: B449=65068, E449=13.97, F449=63, G449=828.45
$ws.Range("B449").Value = 65068
$ws.Range("E449").Value = 13.97
$ws.Range("F449").Value = 63
$ws.Range("G449").Value = 828.45

# Row 456: B456=45706, E456=23.58, F456=-202, G456=-3985.46
$ws.Range("B456").Value = 45706
$ws.Range("E456").Value = 23.58
$ws.Range("F456").Value = -202
$ws.Range("G456").Value = -3985.46

# Row 457: B457=64922, E457=20.98, F457=67, G457=1321.91
$ws.Range("B457").Value = 64922
$ws.Range("E457").Value = 20.98
$ws.Range("F457").Value = 67
$ws.Range("G457").Value = 1321.91

# Row 563: B563=60025, E563=37.22, F563=-98, G563=-3217.34
$ws.Range("B563").Value = 60025
$ws.Range("E563").Value = 37.22
$ws.Range("F563").Value = -98
$ws.Range("G563").Value = -3217.34

# Row 564: B564=64833, E564=34.9, F564=95, G564=3118.85
$ws.Range("B564").Value = 64833
$ws.Range("E564").Value = 34.9
$ws.Range("F564").Value = 95
$ws.Range("G564").Value = 3118.85

# Row 573: B573=64830, E573=34.9, F573=101, G573=3315.83
$ws.Range("B573").Value = 64830
$ws.Range("E573").Value = 34.9
$ws.Range("F573").Value = 101
$ws.Range("G573").Value = 3315.83

# Row 574: B574=60022, E574=37.22, F574=-113, G574=-3709.79
$ws.Range("B574").Value = 60022
$ws.Range("E574").Value = 37.22
$ws.Range("F574").Value = -113
$ws.Range("G574").Value = -3709.79

# Row 625: F625=79, G625=6548.31
$ws.Range("F625").Value = 79
$ws.Range("G625").Value = 6548.31

# Row 626: F626=534, G626=65067.9
$ws.Range("F626").Value = 534
$ws.Range("G626").Value = 65067.9

# Row 628: F628=84, G628=7187.04
$ws.Range("F628").Value = 84
$ws.Range("G628").Value = 7187.04

# Row 629: F629=86, G629=7358.16
$ws.Range("F629").Value = 86
$ws.Range("G629").Value = 7358.16

# Row 630: F630=78, G630=6465.42
$ws.Range("F630").Value = 78
$ws.Range("G630").Value = 6465.42

# Row 632: F632=83, G632=6879.87
$ws.Range("F632").Value = 83
$ws.Range("G632").Value = 6879.87

# Row 633: B633=118089.74
$ws.Range("B633").Value = 118089.74

# Row 675: F675=59, G675=14482.73
$ws.Range("F675").Value = 59
$ws.Range("G675").Value = 14482.73

# Row 678: B678=274730.92
$ws.Range("B678").Value = 274730.92

# Row 721: B721=65079, F721=6, G721=245.22
$ws.Range("B721").Value = 65079
$ws.Range("F721").Value = 6
$ws.Range("G721").Value = 245.22

# Row 722: B722=65362, F722=18, G722=735.66
$ws.Range("B722").Value = 65362
$ws.Range("F722").Value = 18
$ws.Range("G722").Value = 735.66

# Row 751: F751=3480, G751=567622.8
$ws.Range("F751").Value = 3480
$ws.Range("G751").Value = 567622.8

# Row 754: F754=264, G754=38187.6
$ws.Range("F754").Value = 264
$ws.Range("G754").Value = 38187.6

# Row 758: B758=693782.41
$ws.Range("B758").Value = 693782.41

# Row 771: F771=113, G771=17985.08
$ws.Range("F771").Value = 113
$ws.Range("G771").Value = 17985.08

# Row 776: B776=38357.58
$ws.Range("B776").Value = 38357.58

# Row 777: B777=5414147.25
$ws.Range("B777").Value = 5414147.25

# Row 778: B778=5414147.25
$ws.Range("B778").Value = 5414147.25
